# Append the next day's COVID-19 data point (10 June 2020) as a new
# table row (row 92) under the existing "Tabela1" table on the
# "Covid-19 podatki" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Covid-19 podatki")

$lastRow = 91
$newRow = 92

# 1) Duplicate the previous last data row (91) down into the new row 92.
#    Copy+Insert carries over both the data and the cell formatting of
#    row 91 so the new row starts out as a normal (non-highlighted) data
#    row, consistent with the rest of the table, and the used range /
#    dimension grow automatically.
$ws.Range("A" + $lastRow + ":J" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":J" + $newRow).Insert(-4121)

# 2) Give the newest row the "latest update" highlight (light fill +
#    border) that the sheet uses to mark the most recent entry, by
#    copying the formatting from an already-highlighted row.
$ws.Range("A72:J72").Copy()
$ws.Range("A" + $newRow + ":J" + $newRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Write this day's actual figures into row 92.
$ws.Range("A" + $newRow).Value = 43992
$ws.Range("B" + $newRow).Value = 85626
$ws.Range("C" + $newRow).Value = 758
$ws.Range("D" + $newRow).Value = 1488
$ws.Range("E" + $newRow).Value = 0
$ws.Range("F" + $newRow).Value = 6
$ws.Range("G" + $newRow).Value = 0
$ws.Range("H" + $newRow).Value = 0
$ws.Range("I" + $newRow).Value = 109
$ws.Range("J" + $newRow).Value = 0

# 4) Grow the table ("Tabela1") / its AutoFilter so the new row becomes
#    part of the table, not just loose data below it.
$lo = $ws.ListObjects.Item("Tabela1")
$lo.Resize($ws.Range("A1:J" + $newRow))

# 5) Match the selection the author ended up with after adding the row.
$ws.Range("A" + $newRow + ":J" + $newRow).Select()
